{"js": "// The document body contains a single table of 20 rows x 5 columns, each\n// cell holding a short arithmetic expression (e.g. \"26+62=\"). The commit\n// replaces the text of every cell with a new expression; the table shape\n// and all run/paragraph formatting (fonts, sizes, alignment) stay the same.\nconst newValues = [\n  [\"82-36=\", \"82-54=\", \"51+9=\", \"65-17=\", \"54-9=\"],\n  [\"56+1=\", \"82+13=\", \"75-50=\", \"27-3=\", \"97-34=\"],\n  [\"95-38=\", \"92-46=\", \"39+54=\", \"56-15=\", \"5-3=\"],\n  [\"3+24=\", \"47+8=\", \"31+56=\", \"71-6=\", \"3+9=\"],\n  [\"57-55=\", \"89-65=\", \"55+9=\", \"50+23=\", \"90+1=\"],\n  [\"58-26=\", \"5+85=\", \"68-4=\", \"83-31=\", \"93-74=\"],\n  [\"94-72=\", \"44+19=\", \"61+3=\", \"27+42=\", \"41-11=\"],\n  [\"52-3=\", \"20+32=\", \"54-49=\", \"40+56=\", \"38-9=\"],\n  [\"70-67=\", \"67-39=\", \"89-74=\", \"26+51=\", \"56+42=\"],\n  [\"86-81=\", \"50+34=\", \"27+59=\", \"55+3=\", \"63+8=\"],\n  [\"5+3=\", \"1+17=\", \"67-41=\", \"49-29=\", \"25+66=\"],\n  [\"30+35=\", \"66+14=\", \"84-39=\", \"52+6=\", \"9+47=\"],\n  [\"43+21=\", \"32-3=\", \"53+30=\", \"37+23=\", \"15+0=\"],\n  [\"75-12=\", \"68+0=\", \"62-20=\", \"53-34=\", \"56+38=\"],\n  [\"45+52=\", \"98-22=\", \"5+13=\", \"78-32=\", \"26+73=\"],\n  [\"67+9=\", \"7+22=\", \"36-0=\", \"13+45=\", \"25+59=\"],\n  [\"9+78=\", \"19+56=\", \"16+5=\", \"88-64=\", \"32-26=\"],\n  [\"41-12=\", \"61-41=\", \"26-0=\", \"44-7=\", \"74+22=\"],\n  [\"18+38=\", \"39+29=\", \"91-72=\", \"11+28=\", \"43-34=\"],\n  [\"63-30=\", \"95-65=\", \"18+22=\", \"29-0=\", \"70-65=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nif (table.rowCount === newValues.length) {\n  // Fast path: replace the whole grid in one shot.\n  table.values = newValues;\n} else {\n  // Defensive fallback in case the table shape differs from what we\n  // expect: patch cell by cell using the min of rows/cols available.\n  const rows = Math.min(table.rowCount, newValues.length);\n  for (let r = 0; r < rows; r++) {\n    const cols = Math.min(table.values[r].length, newValues[r].length);\n    for (let c = 0; c < cols; c++) {\n      const cell = table.getCell(r, c);\n      cell.body.clear();\n      cell.body.insertText(newValues[r][c], Word.InsertLocation.start);\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document body contains a single table of 20 rows x 5 columns, each\n# cell holding a short arithmetic expression (e.g. \"26+62=\"). The commit\n# replaces the text of every cell with a new expression; the table shape\n# and all run/paragraph formatting (fonts, sizes, alignment) stay the same.\n$newValues = @(\n    @(\"82-36=\", \"82-54=\", \"51+9=\", \"65-17=\", \"54-9=\"),\n    @(\"56+1=\", \"82+13=\", \"75-50=\", \"27-3=\", \"97-34=\"),\n    @(\"95-38=\", \"92-46=\", \"39+54=\", \"56-15=\", \"5-3=\"),\n    @(\"3+24=\", \"47+8=\", \"31+56=\", \"71-6=\", \"3+9=\"),\n    @(\"57-55=\", \"89-65=\", \"55+9=\", \"50+23=\", \"90+1=\"),\n    @(\"58-26=\", \"5+85=\", \"68-4=\", \"83-31=\", \"93-74=\"),\n    @(\"94-72=\", \"44+19=\", \"61+3=\", \"27+42=\", \"41-11=\"),\n    @(\"52-3=\", \"20+32=\", \"54-49=\", \"40+56=\", \"38-9=\"),\n    @(\"70-67=\", \"67-39=\", \"89-74=\", \"26+51=\", \"56+42=\"),\n    @(\"86-81=\", \"50+34=\", \"27+59=\", \"55+3=\", \"63+8=\"),\n    @(\"5+3=\", \"1+17=\", \"67-41=\", \"49-29=\", \"25+66=\"),\n    @(\"30+35=\", \"66+14=\", \"84-39=\", \"52+6=\", \"9+47=\"),\n    @(\"43+21=\", \"32-3=\", \"53+30=\", \"37+23=\", \"15+0=\"),\n    @(\"75-12=\", \"68+0=\", \"62-20=\", \"53-34=\", \"56+38=\"),\n    @(\"45+52=\", \"98-22=\", \"5+13=\", \"78-32=\", \"26+73=\"),\n    @(\"67+9=\", \"7+22=\", \"36-0=\", \"13+45=\", \"25+59=\"),\n    @(\"9+78=\", \"19+56=\", \"16+5=\", \"88-64=\", \"32-26=\"),\n    @(\"41-12=\", \"61-41=\", \"26-0=\", \"44-7=\", \"74+22=\"),\n    @(\"18+38=\", \"39+29=\", \"91-72=\", \"11+28=\", \"43-34=\"),\n    @(\"63-30=\", \"95-65=\", \"18+22=\", \"29-0=\", \"70-65=\"),\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    if ($r -gt $newValues.Length) { break }\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($c -gt $rowValues.Length) { break }\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
